# Updates Betfair Back/Lay odds cells on Sheet1 (data rows 2-14, columns F:AO)
# to match the refreshed odds snapshot described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 2.2  # F2
$ws.Cells.Item(2, 7).Value = 5.5  # G2
$ws.Cells.Item(2, 8).Value = 2.24  # H2
$ws.Cells.Item(2, 9).Value = 990  # I2
$ws.Cells.Item(2, 10).Value = 1.03  # J2
$ws.Cells.Item(2, 11).Value = 9.800000000000001  # K2
$ws.Cells.Item(2, 12).Value = 1.01  # L2
$ws.Cells.Item(2, 13).Value = 1.07  # M2
$ws.Cells.Item(2, 14).Value = 2.72  # N2
$ws.Cells.Item(2, 15).Value = 1.35  # O2
$ws.Cells.Item(2, 16).Value = 1.25  # P2
$ws.Cells.Item(2, 17).Value = 1.37  # Q2
$ws.Cells.Item(2, 18).Value = 1.23  # R2
$ws.Cells.Item(2, 19).Value = 2.2  # S2
$ws.Cells.Item(2, 20).Value = 1.05  # T2
$ws.Cells.Item(2, 21).Value = 1.05  # U2
$ws.Cells.Item(2, 22).Value = 1.53  # V2
$ws.Cells.Item(2, 23).Value = 1.33  # W2
$ws.Cells.Item(2, 27).Value = 1000  # AA2
$ws.Cells.Item(2, 28).Value = 1000  # AB2
$ws.Cells.Item(2, 29).Value = 1000  # AC2
$ws.Cells.Item(2, 30).Value = 1000  # AD2
$ws.Cells.Item(2, 31).Value = 1000  # AE2
$ws.Cells.Item(2, 38).Value = 1000  # AL2
$ws.Cells.Item(2, 39).Value = 1000  # AM2

# Row 3
$ws.Cells.Item(3, 6).Value = 2.54  # F3
$ws.Cells.Item(3, 7).Value = 3.15  # G3
$ws.Cells.Item(3, 8).Value = 2.58  # H3
$ws.Cells.Item(3, 9).Value = 3.25  # I3
$ws.Cells.Item(3, 10).Value = 3.25  # J3
$ws.Cells.Item(3, 11).Value = 4.2  # K3
$ws.Cells.Item(3, 12).Value = 1.01  # L3
$ws.Cells.Item(3, 13).Value = 1.04  # M3
$ws.Cells.Item(3, 14).Value = 1.25  # N3
$ws.Cells.Item(3, 15).Value = 1.27  # O3
$ws.Cells.Item(3, 16).Value = 1.85  # P3
$ws.Cells.Item(3, 17).Value = 1.73  # Q3
$ws.Cells.Item(3, 18).Value = 1.18  # R3
$ws.Cells.Item(3, 19).Value = 1.27  # S3
$ws.Cells.Item(3, 20).Value = 1.05  # T3
$ws.Cells.Item(3, 21).Value = 1.05  # U3
$ws.Cells.Item(3, 22).Value = 1.45  # V3
$ws.Cells.Item(3, 23).Value = 1.56  # W3
$ws.Cells.Item(3, 24).Value = 1000  # X3
$ws.Cells.Item(3, 25).Value = 1000  # Y3
$ws.Cells.Item(3, 26).Value = 1000  # Z3
$ws.Cells.Item(3, 28).Value = 1000  # AB3
$ws.Cells.Item(3, 29).Value = 1000  # AC3
$ws.Cells.Item(3, 30).Value = 1000  # AD3
$ws.Cells.Item(3, 32).Value = 1000  # AF3
$ws.Cells.Item(3, 33).Value = 1000  # AG3
$ws.Cells.Item(3, 34).Value = 1000  # AH3
$ws.Cells.Item(3, 40).Value = 1000  # AN3
$ws.Cells.Item(3, 41).Value = 1000  # AO3

# Row 4
$ws.Cells.Item(4, 6).Value = 2.34  # F4
$ws.Cells.Item(4, 7).Value = 2.64  # G4
$ws.Cells.Item(4, 8).Value = 2.92  # H4
$ws.Cells.Item(4, 9).Value = 3.75  # I4
$ws.Cells.Item(4, 10).Value = 3  # J4
$ws.Cells.Item(4, 11).Value = 3.9  # K4
$ws.Cells.Item(4, 12).Value = 1.01  # L4
$ws.Cells.Item(4, 13).Value = 1.05  # M4
$ws.Cells.Item(4, 14).Value = 2.86  # N4
$ws.Cells.Item(4, 15).Value = 1.29  # O4
$ws.Cells.Item(4, 16).Value = 1.89  # P4
$ws.Cells.Item(4, 18).Value = 1.29  # R4
$ws.Cells.Item(4, 19).Value = 2.32  # S4
$ws.Cells.Item(4, 20).Value = 1.05  # T4
$ws.Cells.Item(4, 21).Value = 1.05  # U4
$ws.Cells.Item(4, 22).Value = 1.36  # V4
$ws.Cells.Item(4, 23).Value = 1.61  # W4
$ws.Cells.Item(4, 24).Value = 1000  # X4
$ws.Cells.Item(4, 28).Value = 15.5  # AB4
$ws.Cells.Item(4, 29).Value = 11.5  # AC4
$ws.Cells.Item(4, 31).Value = 1000  # AE4

# Row 5
$ws.Cells.Item(5, 6).Value = 8.800000000000001  # F5
$ws.Cells.Item(5, 7).Value = 9  # G5
$ws.Cells.Item(5, 8).Value = 1.42  # H5
$ws.Cells.Item(5, 9).Value = 1.43  # I5
$ws.Cells.Item(5, 10).Value = 5.4  # J5
$ws.Cells.Item(5, 14).Value = 5.4  # N5
$ws.Cells.Item(5, 15).Value = 1.21  # O5
$ws.Cells.Item(5, 16).Value = 2.54  # P5
$ws.Cells.Item(5, 17).Value = 1.63  # Q5
$ws.Cells.Item(5, 18).Value = 1.58  # R5
$ws.Cells.Item(5, 19).Value = 2.64  # S5
$ws.Cells.Item(5, 20).Value = 1.84  # T5
$ws.Cells.Item(5, 21).Value = 2.06  # U5
$ws.Cells.Item(5, 22).Value = 3.3  # V5
$ws.Cells.Item(5, 25).Value = 10  # Y5
$ws.Cells.Item(5, 27).Value = 12  # AA5
$ws.Cells.Item(5, 29).Value = 12  # AC5
$ws.Cells.Item(5, 30).Value = 9.800000000000001  # AD5
$ws.Cells.Item(5, 32).Value = 80  # AF5
$ws.Cells.Item(5, 33).Value = 32  # AG5
$ws.Cells.Item(5, 38).Value = 100  # AL5
$ws.Cells.Item(5, 41).Value = 5.5  # AO5

# Row 6
$ws.Cells.Item(6, 8).Value = 1.72  # H6
$ws.Cells.Item(6, 12).Value = 1.01  # L6
$ws.Cells.Item(6, 14).Value = 5  # N6
$ws.Cells.Item(6, 15).Value = 1.16  # O6
$ws.Cells.Item(6, 17).Value = 1.46  # Q6
$ws.Cells.Item(6, 18).Value = 1.24  # R6
$ws.Cells.Item(6, 19).Value = 1.72  # S6
$ws.Cells.Item(6, 20).Value = 1.05  # T6
$ws.Cells.Item(6, 21).Value = 1.05  # U6
$ws.Cells.Item(6, 23).Value = 1.18  # W6
$ws.Cells.Item(6, 24).Value = 1000  # X6
$ws.Cells.Item(6, 25).Value = 1000  # Y6
$ws.Cells.Item(6, 26).Value = 1000  # Z6
$ws.Cells.Item(6, 28).Value = 1000  # AB6
$ws.Cells.Item(6, 29).Value = 1000  # AC6
$ws.Cells.Item(6, 30).Value = 1000  # AD6
$ws.Cells.Item(6, 31).Value = 1000  # AE6
$ws.Cells.Item(6, 32).Value = 1000  # AF6
$ws.Cells.Item(6, 33).Value = 1000  # AG6
$ws.Cells.Item(6, 34).Value = 1000  # AH6
$ws.Cells.Item(6, 35).Value = 1000  # AI6
$ws.Cells.Item(6, 36).Value = 1000  # AJ6
$ws.Cells.Item(6, 37).Value = 1000  # AK6
$ws.Cells.Item(6, 38).Value = 1000  # AL6
$ws.Cells.Item(6, 39).Value = 1000  # AM6
$ws.Cells.Item(6, 40).Value = 1000  # AN6
$ws.Cells.Item(6, 41).Value = 1000  # AO6

# Row 7
$ws.Cells.Item(7, 6).Value = 4.3  # F7
$ws.Cells.Item(7, 17).Value = 1.77  # Q7
$ws.Cells.Item(7, 19).Value = 2.98  # S7
$ws.Cells.Item(7, 26).Value = 12.5  # Z7
$ws.Cells.Item(7, 27).Value = 21  # AA7
$ws.Cells.Item(7, 28).Value = 19.5  # AB7
$ws.Cells.Item(7, 41).Value = 10.5  # AO7

# Row 8
$ws.Cells.Item(8, 10).Value = 1.03  # J8
$ws.Cells.Item(8, 13).Value = 1.06  # M8
$ws.Cells.Item(8, 14).Value = 1.1  # N8
$ws.Cells.Item(8, 16).Value = 1.25  # P8
$ws.Cells.Item(8, 17).Value = 1.01  # Q8
$ws.Cells.Item(8, 18).Value = 1.38  # R8
$ws.Cells.Item(8, 19).Value = 1.71  # S8
$ws.Cells.Item(8, 22).Value = 1.02  # V8
$ws.Cells.Item(8, 23).Value = 4.1  # W8
$ws.Cells.Item(8, 33).Value = 17  # AG8

# Row 9
$ws.Cells.Item(9, 8).Value = 3.5  # H9
$ws.Cells.Item(9, 11).Value = 3.45  # K9
$ws.Cells.Item(9, 12).Value = 1.46  # L9
$ws.Cells.Item(9, 14).Value = 3.7  # N9
$ws.Cells.Item(9, 16).Value = 1.87  # P9
$ws.Cells.Item(9, 17).Value = 2.12  # Q9
$ws.Cells.Item(9, 18).Value = 1.33  # R9
$ws.Cells.Item(9, 19).Value = 3.85  # S9
$ws.Cells.Item(9, 20).Value = 1.84  # T9
$ws.Cells.Item(9, 21).Value = 2.12  # U9
$ws.Cells.Item(9, 22).Value = 1.39  # V9
$ws.Cells.Item(9, 24).Value = 12.5  # X9
$ws.Cells.Item(9, 28).Value = 9.6  # AB9
$ws.Cells.Item(9, 31).Value = 42  # AE9
$ws.Cells.Item(9, 34).Value = 18  # AH9
$ws.Cells.Item(9, 38).Value = 42  # AL9
$ws.Cells.Item(9, 39).Value = 110  # AM9
$ws.Cells.Item(9, 40).Value = 21  # AN9
$ws.Cells.Item(9, 41).Value = 42  # AO9

# Row 10
$ws.Cells.Item(10, 7).Value = 14  # G10
$ws.Cells.Item(10, 17).Value = 1.74  # Q10
$ws.Cells.Item(10, 21).Value = 1.74  # U10
$ws.Cells.Item(10, 23).Value = 1.07  # W10
$ws.Cells.Item(10, 28).Value = 40  # AB10
$ws.Cells.Item(10, 29).Value = 13.5  # AC10
$ws.Cells.Item(10, 31).Value = 14.5  # AE10
$ws.Cells.Item(10, 39).Value = 220  # AM10

# Row 11
$ws.Cells.Item(11, 6).Value = 2.74  # F11
$ws.Cells.Item(11, 7).Value = 2.76  # G11
$ws.Cells.Item(11, 8).Value = 2.84  # H11
$ws.Cells.Item(11, 9).Value = 2.86  # I11
$ws.Cells.Item(11, 10).Value = 3.5  # J11
$ws.Cells.Item(11, 11).Value = 3.55  # K11
$ws.Cells.Item(11, 12).Value = 1.41  # L11
$ws.Cells.Item(11, 22).Value = 1.54  # V11
$ws.Cells.Item(11, 23).Value = 1.56  # W11
$ws.Cells.Item(11, 26).Value = 18.5  # Z11
$ws.Cells.Item(11, 29).Value = 7.6  # AC11
$ws.Cells.Item(11, 31).Value = 30  # AE11
$ws.Cells.Item(11, 37).Value = 29  # AK11

# Row 12
$ws.Cells.Item(12, 7).Value = 4.3  # G12
$ws.Cells.Item(12, 9).Value = 1.95  # I12
$ws.Cells.Item(12, 17).Value = 1.69  # Q12
$ws.Cells.Item(12, 19).Value = 2.74  # S12
$ws.Cells.Item(12, 20).Value = 1.61  # T12
$ws.Cells.Item(12, 22).Value = 2.04  # V12
$ws.Cells.Item(12, 24).Value = 20  # X12
$ws.Cells.Item(12, 25).Value = 11.5  # Y12
$ws.Cells.Item(12, 28).Value = 19.5  # AB12
$ws.Cells.Item(12, 36).Value = 80  # AJ12
$ws.Cells.Item(12, 41).Value = 10  # AO12

# Row 13
$ws.Cells.Item(13, 6).Value = 4.8  # F13
$ws.Cells.Item(13, 7).Value = 4.9  # G13
$ws.Cells.Item(13, 8).Value = 1.79  # H13
$ws.Cells.Item(13, 9).Value = 1.8  # I13
$ws.Cells.Item(13, 17).Value = 1.58  # Q13
$ws.Cells.Item(13, 20).Value = 1.6  # T13
$ws.Cells.Item(13, 21).Value = 2.56  # U13
$ws.Cells.Item(13, 22).Value = 2.24  # V13
$ws.Cells.Item(13, 23).Value = 1.25  # W13
$ws.Cells.Item(13, 30).Value = 10.5  # AD13

# Row 14
$ws.Cells.Item(14, 6).Value = 1.32  # F14
$ws.Cells.Item(14, 8).Value = 1.97  # H14
$ws.Cells.Item(14, 9).Value = 980  # I14
$ws.Cells.Item(14, 10).Value = 3.3  # J14
$ws.Cells.Item(14, 11).Value = 980  # K14
$ws.Cells.Item(14, 12).Value = 1.01  # L14
$ws.Cells.Item(14, 13).Value = 1.09  # M14
$ws.Cells.Item(14, 14).Value = 1.63  # N14
$ws.Cells.Item(14, 15).Value = 1.41  # O14
$ws.Cells.Item(14, 16).Value = 1.25  # P14
$ws.Cells.Item(14, 17).Value = 1.01  # Q14
$ws.Cells.Item(14, 18).Value = 1.13  # R14
$ws.Cells.Item(14, 19).Value = 2.3  # S14
$ws.Cells.Item(14, 20).Value = 1.05  # T14
$ws.Cells.Item(14, 21).Value = 1.05  # U14
$ws.Cells.Item(14, 25).Value = 1000  # Y14
$ws.Cells.Item(14, 27).Value = 1000  # AA14
$ws.Cells.Item(14, 28).Value = 1000  # AB14
$ws.Cells.Item(14, 29).Value = 1000  # AC14
$ws.Cells.Item(14, 30).Value = 1000  # AD14
$ws.Cells.Item(14, 32).Value = 1000  # AF14
$ws.Cells.Item(14, 33).Value = 1000  # AG14
$ws.Cells.Item(14, 39).Value = 1000  # AM14
$ws.Cells.Item(14, 40).Value = 1000  # AN14
